# Auto-generated edit script: updates crypto price/volume table
# to match the refreshed data snapshot from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.984.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.058.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.058.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.582.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000159"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.024.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.055.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "346.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.495"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0838"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.21%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.73%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0651"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.686"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.397.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.096.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0259"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.921"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.85%  "
